$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.434.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.554.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.92%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.554.23"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.96%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("E10").Value = "  +2.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.06"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.394"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.162.28"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000188"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.07"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.555.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.22%  "
$ws.Range("E17").Value = "  +1.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.523.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("E21").Value = "  +3.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "395.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.572"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.18%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.701.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.89%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  +8.67%  "
$ws.Range("E28").Value = "  +7.74%  "
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.575.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.91%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  +1.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "168.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0802"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.68%  "
$ws.Range("E42").Value = "  +1.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +14.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.65%  "
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.54%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.13%  "
$ws.Range("E48").Value = "  +6.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.433.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.72%  "
